$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(
    "9/1/2012","9/2/2012","9/3/2012","9/4/2012","9/5/2012","9/6/2012","9/7/2012",
    "9/8/2012","9/9/2012","9/10/2012","9/11/2012","9/12/2012","9/13/2012","9/14/2012",
    "9/15/2012","9/16/2012","9/17/2012","9/18/2012","9/19/2012","9/20/2012","9/21/2012",
    "9/22/2012","9/23/2012","9/24/2012","9/25/2012","9/26/2012","9/27/2012","9/28/2012",
    "9/29/2012","9/30/2012","10/1/2012","10/2/2012"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}
